$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp string (14:22 -> 14:52)
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 14:52"

# 2) Arabia Saudita's updated case count (4934) now overtakes Filipinas
#    (4932) in the ranking, so it moves up to row 37 - right after
#    Pakistan - and the countries that used to occupy rows 37-40
#    (Filipinas, Malasia, Mexico, Indonesia) each shift down one row.
#    Emiratos Arabes Unidos (row 42) and every other row is unaffected.

$ws.Range("A37").Value = "Arabia Saudita"
$ws.Range("B37").Value = 4934
$ws.Range("C37").Value = 472
$ws.Range("D37").Value = 805
$ws.Range("E37").Value = 4064
$ws.Range("F37").Value = 67
$ws.Range("G37").Value = 6
$ws.Range("H37").Value = 65

$ws.Range("A38").Value = "Filipinas"
$ws.Range("B38").Value = 4932
$ws.Range("C38").Value = 284
$ws.Range("D38").Value = 242
$ws.Range("E38").Value = 4375
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = 18
$ws.Range("H38").Value = 315

$ws.Range("A39").Value = "Malasia"
$ws.Range("B39").Value = 4817
$ws.Range("C39").Value = 134
$ws.Range("D39").Value = 2276
$ws.Range("E39").Value = 2464
$ws.Range("F39").Value = 66
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 77

$ws.Range("A40").Value = "Mexico"
$ws.Range("B40").Value = 4661
$ws.Range("C40").Value = 442
$ws.Range("D40").Value = 1843
$ws.Range("E40").Value = 2522
$ws.Range("F40").Value = 185
$ws.Range("G40").Value = 23
$ws.Range("H40").Value = 296

$ws.Range("A41").Value = "Indonesia"
$ws.Range("B41").Value = 4557
$ws.Range("C41").Value = 316
$ws.Range("D41").Value = 380
$ws.Range("E41").Value = 3778
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 26
$ws.Range("H41").Value = 399
